$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing molecule value and append the two new molecules
$ws.Range("A2").Value = "Adrenaline (1Mg)"
$ws.Range("A3").Value = "Acyclovir (800Mg)"
$ws.Range("A4").Value = "Acetylcysteine-200Mg"

# Move active selection to the next empty row, as in the saved workbook
$ws.Range("A5").Select()
